# Wakanda interface selectors - "bundle" column update
#   1. Renumber existing vpc00X (3-digit) values in rows 3..62 to vpc000X (4-digit).
#   2. Extend the table with 48 new rows (63..110) continuing the same
#      zero-padded sequence (vpc0061..vpc0108), copying the existing
#      alternating row styling from the last two data rows (61/62).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("formulas")

# --- 1. Fix up the already-present entries (rows 3-62, C column) ---------
for ($row = 3; $row -le 62; $row++) {
    $num = "{0:D4}" -f ($row - 2)
    $ws.Cells.Item($row, 3).Value = "vpc" + $num
}

# --- 2. Grow the table: copy the style of the last odd/even row pair -----
#        (rows 61-62) down across the 48 new rows (63-110) ---------------
$ws.Range("A61:D62").Copy()
$ws.Range("A63:D110").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 3. Populate the new bundle values (vpc0061 .. vpc0108) --------------
for ($row = 63; $row -le 110; $row++) {
    $num = "{0:D4}" -f ($row - 2)
    $ws.Cells.Item($row, 3).Value = "vpc" + $num
}
